# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the per-class Leve profit sheets (H:N columns) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 259.64285
$ws.Range("J28").Value = 145.5
$ws.Range("L28").Value = 145.5
$ws.Range("N28").Value = -1115.5

$ws.Range("H38").Value = 1141
$ws.Range("I38").Value = 47.4
$ws.Range("J38").Value = 3875
$ws.Range("K38").Value = 142.2
$ws.Range("L38").Value = 11625
$ws.Range("M38").Value = 229.8
$ws.Range("N38").Value = -12369

$ws.Range("H62").Value = 4369.4165
$ws.Range("I62").Value = 3324.5
$ws.Range("K62").Value = 3324.5
$ws.Range("M62").Value = -2700.5

$ws.Range("H65").Value = 4369.4165
$ws.Range("I65").Value = 3324.5
$ws.Range("K65").Value = 16622.5
$ws.Range("M65").Value = -13502.5

$ws.Range("H86").Value = 6326.0586
$ws.Range("I86").Value = 7938.5835
$ws.Range("J86").Value = 2456
$ws.Range("K86").Value = 7938.5835
$ws.Range("L86").Value = 2456
$ws.Range("M86").Value = -6815.5835
$ws.Range("N86").Value = -4702

$ws.Range("H89").Value = 6326.0586
$ws.Range("I89").Value = 7938.5835
$ws.Range("J89").Value = 2456
$ws.Range("K89").Value = 39692.9175
$ws.Range("L89").Value = 12280
$ws.Range("M89").Value = -34076.9175
$ws.Range("N89").Value = -23512

$ws.Range("H111").Value = 532.3
$ws.Range("J111").Value = 425
$ws.Range("L111").Value = 1275
$ws.Range("N111").Value = -7409

$ws.Range("H118").Value = 655.8570999999999
$ws.Range("I118").Value = 619.5
$ws.Range("J118").Value = 874
$ws.Range("K118").Value = 1858.5
$ws.Range("L118").Value = 2622
$ws.Range("M118").Value = -201.5
$ws.Range("N118").Value = -5936

$ws.Range("H125").Value = 6355
$ws.Range("I125").Value = 2862.5
$ws.Range("J125").Value = 11011.667
$ws.Range("K125").Value = 25762.5
$ws.Range("L125").Value = 99105.003
$ws.Range("M125").Value = -23302.5
$ws.Range("N125").Value = -104025.003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11443.1
$ws.Range("I32").Value = 12581.223
$ws.Range("J32").Value = 1200
$ws.Range("K32").Value = 12581.223
$ws.Range("L32").Value = 1200
$ws.Range("M32").Value = -12294.223
$ws.Range("N32").Value = -1774

$ws.Range("H45").Value = 1836.1818
$ws.Range("J45").Value = 1868
$ws.Range("L45").Value = 1868
$ws.Range("N45").Value = -2622

$ws.Range("H132").Value = 3504.5
$ws.Range("I132").Value = 3504.5
$ws.Range("K132").Value = 10513.5
$ws.Range("M132").Value = -7983.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2942.9412
$ws.Range("I86").Value = 2247.1
$ws.Range("J86").Value = 3937
$ws.Range("K86").Value = 2247.1
$ws.Range("L86").Value = 3937
$ws.Range("M86").Value = -1124.1
$ws.Range("N86").Value = -6183

$ws.Range("H89").Value = 2942.9412
$ws.Range("I89").Value = 2247.1
$ws.Range("J89").Value = 3937
$ws.Range("K89").Value = 11235.5
$ws.Range("L89").Value = 19685
$ws.Range("M89").Value = -5619.5
$ws.Range("N89").Value = -30917

$ws.Range("H99").Value = 1642.5714
$ws.Range("I99").Value = 1649.6666
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 1649.6666
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = -151.6666
$ws.Range("N99").Value = -4596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8999.143
$ws.Range("J134").Value = 9333
$ws.Range("L134").Value = 27999
$ws.Range("N134").Value = -33069

$ws.Range("H135").Value = 48399
$ws.Range("J135").Value = 48399
$ws.Range("L135").Value = 48399
$ws.Range("N135").Value = -58539

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1500
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5685

$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -4908

$ws.Range("H137").Value = 2851
$ws.Range("J137").Value = 3244
$ws.Range("L137").Value = 9732
$ws.Range("N137").Value = -19932

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3099.7144
$ws.Range("I80").Value = 2274.75
$ws.Range("J80").Value = 4199.6665
$ws.Range("K80").Value = 2274.75
$ws.Range("L80").Value = 4199.6665
$ws.Range("M80").Value = -1276.75
$ws.Range("N80").Value = -6195.6665

$ws.Range("H83").Value = 3099.7144
$ws.Range("I83").Value = 2274.75
$ws.Range("J83").Value = 4199.6665
$ws.Range("K83").Value = 11373.75
$ws.Range("L83").Value = 20998.3325
$ws.Range("M83").Value = -6381.75
$ws.Range("N83").Value = -30982.3325

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()

$ws.Range("H113").Value = 1800
$ws.Range("I113").Value = 1800
$ws.Range("K113").Value = 1800
$ws.Range("M113").Value = 370

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 649.25
$ws.Range("I9").Value = 649.25
$ws.Range("K9").Value = 649.25
$ws.Range("M9").Value = -425.25

$ws.Range("H16").Value = 2679.1
$ws.Range("I16").Value = 1435.8572
$ws.Range("K16").Value = 1435.8572
$ws.Range("M16").Value = -1265.8572

$ws.Range("H30").Value = 505.77777
$ws.Range("I30").Value = 362
$ws.Range("J30").Value = 1009
$ws.Range("K30").Value = 362
$ws.Range("L30").Value = 1009
$ws.Range("M30").Value = -254
$ws.Range("N30").Value = -1225

$ws.Range("H55").Value = 233.22223
$ws.Range("I55").Value = 216.5
$ws.Range("J55").Value = 266.66666
$ws.Range("K55").Value = 216.5
$ws.Range("L55").Value = 266.66666
$ws.Range("M55").Value = -43.5
$ws.Range("N55").Value = -612.66666

$ws.Range("H122").Value = 3710.3635
$ws.Range("I122").Value = 3336
$ws.Range("J122").Value = 4159.6
$ws.Range("K122").Value = 10008
$ws.Range("L122").Value = 12478.8
$ws.Range("M122").Value = -7558
$ws.Range("N122").Value = -17378.8

$ws.Range("H136").Value = 3347.182
$ws.Range("J136").Value = 2005
$ws.Range("L136").Value = 6015
$ws.Range("N136").Value = -11115

$ws.Range("H138").Value = 49994
$ws.Range("J138").Value = 49994
$ws.Range("L138").Value = 49994
$ws.Range("N138").Value = -60274
